# Updated symbol list (cryptos.xlsx) - price/ranking refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.52"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").Value = "'22.45"
$ws.Range("D3").Style = "Normal"

$ws.Range("D4").Value = "'5.253"
$ws.Range("D4").Style = "Normal"

$ws.Range("D6").Value = "'3.419"
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").Value = "'6.309"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").Value = "'0.8075"
$ws.Range("D8").Style = "Normal"

$ws.Range("D9").Value = "'0.8715"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").Value = "'0.1411"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").Value = "'0.07397"
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").Value = "'0.03040"
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").Value = "'0.03076"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'0.09384"
$ws.Range("D14").Style = "Normal"

$ws.Range("D15").Value = "'3.872"
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").Value = "'0.001581"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").Value = "'0.04770"
$ws.Range("D17").Style = "Normal"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005855"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17OneONE"

$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").Value = "'0.006440"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18TigerCashTCH"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").Value = "'0.005030"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").Value = "'0.0009972"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20BitKanKAN"

$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21NitroExNTX"

$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Value = "'3.691"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22LEOLEO"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Value = "'2.195"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23BTSETokenBTSE"

$ws.Range("D25").Value = "'0.3278"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'0.1283"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").Value = "'0.01830"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "26UpBotsUBXTBestin24h"

$ws.Range("D40").Value = "'0.03930"
$ws.Range("D40").Style = "Normal"

$ws.Range("D41").Value = "'0.006793"
$ws.Range("D41").Style = "Normal"

$ws.Range("D43").Value = "'0.002671"
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").Value = "'0.008459"
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").Value = "'0.00005591"
$ws.Range("D45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'0.4504"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'0.1967"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
